$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 adds "5000.0" (stored as text, like the header row) under both the
# Meta (A) and Venda (B) columns.
#
# Forcing NumberFormat to "@" (Text) before the write keeps Excel from
# reinterpreting the numeric-looking string as a number; resetting the
# style back to "Normal" afterwards drops the scratch text format so the
# cells end up on the default style again, matching the source data.
$range = $ws.Range("A2:B2")
$range.NumberFormat = "@"
$ws.Range("A2").Value = "5000.0"
$ws.Range("B2").Value = "5000.0"
$range.Style = "Normal"
